$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# DATE_TYPE_CODE: "001" -> "002" (keep as text, matching existing General-formatted
# text cell convention used throughout this row)
$ws.Range("J2").Formula = "'002"

# REPORT_DATE: updated reporting period
$ws.Range("N2").Value = "2020-06-30 00:00:00"

# Financial figures updated for the new report date
$ws.Range("O2").Value = 44071294.14
$ws.Range("P2").Value = 488068692.66
$ws.Range("Q2").Value = 428381181.29

# TOE_RATIO cleared (no longer supplied)
$ws.Range("R2").Value = ""

$ws.Range("S2").Value = 387860071.43
$ws.Range("T2").Value = 387860071.43

# OPERATE_EXPENSE_RATIO cleared (no longer supplied)
$ws.Range("U2").Value = ""

$ws.Range("V2").Value = 4710234.38
$ws.Range("W2").Value = 32841329.34
$ws.Range("X2").Value = 54502.32
$ws.Range("Y2").Value = 62546378.87
$ws.Range("Z2").Value = 62482519.65
$ws.Range("AA2").Value = 15260027.13
$ws.Range("AG2").Value = 2915043.82
$ws.Range("AP2").Value = 21.8469833909
$ws.Range("AQ2").Value = 57.055420470867
$ws.Range("AR2").Value = 61.308047011112
$ws.Range("AS2").Value = 42115494.14
$ws.Range("AT2").Value = 34.49544448596
